$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename transportation item label (A29, merged A29:A30) in sharedStrings
$ws.Range("A29").Value = "Trucking [tonne*km]"

# Updated numeric results (Bwaise system B LCA recalculation)
$ws.Range("C2").Value = 885937.5
$ws.Range("E2").Value = 248062.5
$ws.Range("F2").Value = 0.06034129842869993
$ws.Range("C3").Value = 885937.5
$ws.Range("E3").Value = 248062.5
$ws.Range("F3").Value = 0.06034129842869993
$ws.Range("C4").Value = 2734375
$ws.Range("E4").Value = 2953125
$ws.Range("F4").Value = 0.7183487908178563
$ws.Range("C5").Value = 2734375
$ws.Range("E5").Value = 2953125
$ws.Range("F5").Value = 0.7183487908178563
$ws.Range("C10").Value = 14281.25
$ws.Range("D10").Value = 0.9853385079775765
$ws.Range("E10").Value = 7569.062500000001
$ws.Range("F10").Value = 0.001841177361100455
$ws.Range("D11").Value = 0.01466149202242346
$ws.Range("C12").Value = 14493.75
$ws.Range("E12").Value = 7681.687500000001
$ws.Range("F12").Value = 0.001868573435620111
$ws.Range("C13").Value = 5000000
$ws.Range("D13").Value = 0.9778357235984355
$ws.Range("E13").Value = 75000
$ws.Range("F13").Value = 0.01824377881442175
$ws.Range("D14").Value = 0.02216427640156454
$ws.Range("C15").Value = 5113333.333333333
$ws.Range("E15").Value = 76699.99999999999
$ws.Range("F15").Value = 0.01865730446754863
$ws.Range("C16").Value = 39375
$ws.Range("E16").Value = 77568.75
$ws.Range("F16").Value = 0.01886862823881569
$ws.Range("C17").Value = 39375
$ws.Range("E17").Value = 77568.75
$ws.Range("F17").Value = 0.01886862823881569
$ws.Range("C18").Value = 12392187.5
$ws.Range("E18").Value = 148706.25
$ws.Range("F18").Value = 0.03617285244429472
$ws.Range("C19").Value = 12392187.5
$ws.Range("E19").Value = 148706.25
$ws.Range("F19").Value = 0.03617285244429472
$ws.Range("C20").Value = 131152.34375
$ws.Range("E20").Value = 334438.4765625001
$ws.Range("F20").Value = 0.0813522879125123
$ws.Range("C22").Value = 131152.34375
$ws.Range("E22").Value = 334438.4765625001
$ws.Range("F22").Value = 0.0813522879125123
$ws.Range("C23").Value = 742.1875
$ws.Range("E23").Value = 146210.9375
$ws.Range("F23").Value = 0.0355658667199899
$ws.Range("C24").Value = 742.1875
$ws.Range("E24").Value = 146210.9375
$ws.Range("F24").Value = 0.0355658667199899
$ws.Range("C29").Value = 675000.0000000007
$ws.Range("C30").Value = 675000.0000000007
$ws.Range("D35").Value = 1.941708215462239
$ws.Range("C36").Value = 44590.22803355606
$ws.Range("D36").Value = 0.007309787903568453
$ws.Range("B37").Value = 768085.007480899
$ws.Range("C37").Value = -2306774.376275815
$ws.Range("D37").Value = -0.3781553083620297
$ws.Range("D38").Value = -0.06662868275258886
$ws.Range("B39").Value = 418058.5060666984
$ws.Range("C39").Value = -2257515.932760172
$ws.Range("D39").Value = -0.3700802481876728
$ws.Range("D40").Value = -0.09504328350191121
$ws.Range("B42").Value = 6166.182812556733
$ws.Range("C42").Value = -33297.38718780636
$ws.Range("D42").Value = -0.005458524183879381
$ws.Range("C44").Value = 6100071.386720845
